$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44537
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 800
$ws.Range("L2").Value = 900
$ws.Range("M2").Value = 850
$ws.Range("P2").Value = 850

# Row 3
$ws.Range("D3").Value = 44827
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 2200
$ws.Range("L3").Value = 2300
$ws.Range("M3").Value = 2250
$ws.Range("P3").Value = 2250

# Row 4
$ws.Range("D4").Value = 44524
$ws.Range("J4").Value = 400

# Row 5
$ws.Range("D5").Value = 44508
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 900
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = 950
$ws.Range("P5").Value = 950

# Row 6
$ws.Range("D6").Value = 44516
$ws.Range("K6").Value = 900
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = 950
$ws.Range("P6").Value = 950

# Row 7
$ws.Range("D7").Value = 44523
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 800
$ws.Range("L7").Value = 900
$ws.Range("M7").Value = 850
$ws.Range("P7").Value = 850

# Row 8
$ws.Range("D8").Value = 44530
$ws.Range("J8").Value = 300

# Row 10
$ws.Range("D10").Value = 44505
$ws.Range("J10").Value = 440
$ws.Range("K10").Value = 900
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = 950
$ws.Range("P10").Value = 950

# Row 11
$ws.Range("D11").Value = 44532
$ws.Range("J11").Value = 240
$ws.Range("K11").Value = 800
$ws.Range("L11").Value = 900
$ws.Range("M11").Value = 850
$ws.Range("P11").Value = 850

# Row 12
$ws.Range("D12").Value = 44504
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 900
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = 950
$ws.Range("P12").Value = 950

# Row 13
$ws.Range("D13").Value = 44525
$ws.Range("J13").Value = 360

# Row 14
$ws.Range("D14").Value = 44503
$ws.Range("K14").Value = 900
$ws.Range("L14").Value = 1000
$ws.Range("M14").Value = 950
$ws.Range("P14").Value = 950

# Row 15
$ws.Range("D15").Value = 44476
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 1100
$ws.Range("L15").Value = 1200
$ws.Range("M15").Value = 1150
$ws.Range("P15").Value = 1150

# Row 16
$ws.Range("D16").Value = 44512
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("J16").Value = 600
$ws.Range("K16").Value = 900
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = 950
$ws.Range("P16").Value = 950

# Row 17
$ws.Range("D17").Value = 44510
$ws.Range("J17").Value = 600

# Row 18
$ws.Range("D18").Value = 44460
$ws.Range("H18").Value = "Verde"
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = 2200
$ws.Range("L18").Value = 2300
$ws.Range("M18").Value = 2250
$ws.Range("P18").Value = 2250

# Row 19
$ws.Range("D19").Value = 44545
$ws.Range("J19").Value = 4000
$ws.Range("K19").Value = 800
$ws.Range("L19").Value = 900
$ws.Range("M19").Value = 850
$ws.Range("P19").Value = 850

# Row 20
$ws.Range("D20").Value = 44511

# Row 21
$ws.Range("D21").Value = 44517
$ws.Range("J21").Value = 500
$ws.Range("K21").Value = 800
$ws.Range("L21").Value = 900
$ws.Range("M21").Value = 850
$ws.Range("P21").Value = 850

# Row 22
$ws.Range("D22").Value = 44553
$ws.Range("J22").Value = 8000
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = 850
$ws.Range("P22").Value = 850

